$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Förändrad" (C) column date values for rows 2-8 from 45233 (2023-11-03) to 45243 (2023-11-13)
$ws.Range("C2:C8").Value = 45243
